# feat: add 2022-Q1 data
#
# Before: sheets = [2021-Q3, 总计]
# After:  sheets = [2021-Q3, 2022-Q1, 总计]
#   - "2022-Q1" is a brand-new sheet (same column layout as "2021-Q3")
#     holding the fund-holding snapshot for the new quarter.
#   - "总计" (the roll-up sheet) gets a new row inserted at the top of its
#     data (row 2) for "2022-Q1", pushing the existing "2021-Q3" row down.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right after "2021-Q3".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "2022-Q1"

# "总计" got pushed from position 2 -> 3 by the insert above, so fetch a
# fresh handle to it now (stale handles lose their style bookkeeping).
$zongji = $wb.Worksheets.Item(3)

# Headers (same labels as "2021-Q3" apart from column D).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("A2").Value = 0

# Match the header/index-cell formatting used elsewhere in the workbook
# (bold, centered, thin-bordered) by pasting formats from "总计"'s header.
$zongji.Range("B1:D1").Copy()
$newSheet.Range("B1:D1").PasteSpecial(-4122)
$zongji.Range("B1:D1").Copy()
$newSheet.Range("E1:G1").PasteSpecial(-4122)
$zongji.Range("B1").Copy()
$newSheet.Range("H1").PasteSpecial(-4122)
$zongji.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Data row. These look numeric but are stored as plain text in the source
# sheet (keeps the leading zero on the fund code), so force text first.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "000747"
$newSheet.Range("C2").Value = "广发逆向策略灵活配置混合"
$newSheet.Range("D2").Value = "1.25"
$newSheet.Range("E2").Value = "89.61"
$newSheet.Range("F2").Value = "4.03"
$newSheet.Range("G2").Value = "0.0504"
$newSheet.Range("H2").Value = 10

# Undo the stray number-format styling the "@" format left behind so the
# data cells stay plain (unstyled), matching the rest of the workbook.
$zongji.Range("B2:C2").Copy()
$newSheet.Range("B2:C2").PasteSpecial(-4122)
$zongji.Range("C2:D2").Copy()
$newSheet.Range("D2:E2").PasteSpecial(-4122)
$zongji.Range("C2:D2").Copy()
$newSheet.Range("F2:G2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Add a "2022-Q1" row to the top of "总计"'s data, keep "2021-Q3"
#    below it.
# ---------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

$zongji.Range("A3").Value = 1

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 1
$zongji.Range("D2").Value = 0.05

# Re-apply the index-column / data-row formatting (the inserted row starts
# out with a freshly-minted blank style).
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)
$zongji.Range("B3:D3").Copy()
$zongji.Range("B2:D2").PasteSpecial(-4122)

Write-Host "2022-Q1 sheet added; 总计 updated"
